$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, border, centered) from F1 onto the new
# header cells G1:I1 before setting their text.
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

$ws.Range("G1").Value = "MSE_median"
$ws.Range("H1").Value = "MAE_median"
$ws.Range("I1").Value = "Dir_accuracy"

# Row 2 (XGBoostRegressor) - new metric columns
$ws.Range("G2").Value = 0.0005189844392051669
$ws.Range("H2").Value = 0.02278119634655595
$ws.Range("I2").Value = 0.4354609929078014

# Row 3 (Naive) - new metric columns; Dir_accuracy (I3) is left blank, mirroring
# the blank R^2 (F3) cell already present for this row.
$ws.Range("G3").Value = 0.001083194070471167
$ws.Range("H3").Value = 0.03291191380748265

$ws.Range("F3").Copy()
$ws.Range("I3").PasteSpecial(-4122)

Write-Output "applied k-important-features metric columns"
